# Auto update Excel log
# Appends a new log row (row 3) to both the "Proximity" and "Camera" sheets,
# recording an EXIT / Image-Captured event for the Living Room Main Door at
# 2026-01-28 18:45:19 (Hour bucket 18:00).

$wb = $excel.ActiveWorkbook

# --- Proximity sheet: door EXIT event -------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")

# Column A holds a literal "yyyy-mm-dd" style string, not a real date value
# in this log (every other row stores it as text) - prefix with an
# apostrophe so Excel keeps it as literal text instead of auto-converting
# it to a date serial number.
$wsProximity.Range("A3").Value = "'2026-01-28"
$wsProximity.Range("B3").Value = "18:45:19"
$wsProximity.Range("C3").Value = "18:00"
$wsProximity.Range("D3").Value = "Living Room Main Door"
$wsProximity.Range("E3").Value = "EXIT"
$wsProximity.Range("F3").Value = "User EXITED Living Room Main Door"

# --- Camera sheet: image captured event ------------------------------------
$wsCamera = $wb.Worksheets.Item("Camera")

$wsCamera.Range("A3").Value = "'2026-01-28"
$wsCamera.Range("B3").Value = "18:45:19"
$wsCamera.Range("C3").Value = "18:00"
$wsCamera.Range("D3").Value = "Living Room Main Door"
$wsCamera.Range("E3").Value = "Image Captured"
$wsCamera.Range("F3").Value = "Active"
